$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.982.12"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "1.640.00"
$ws.Range("E3").Value = "  -1.86%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'209.79"
$ws.Range("E5").Value = "  -0.49%  "
$ws.Range("D6").Value = "'0.5168"
$ws.Range("E6").Value = "  -1.40%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.2564"
$ws.Range("D9").Value = "'0.06234"
$ws.Range("E9").Value = "  -0.53%  "
$ws.Range("D10").Value = "'20.37"
$ws.Range("E10").Value = "  -3.64%  "
$ws.Range("D11").Value = "'0.07526"
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("D12").Value = "1.641.44"
$ws.Range("E12").Value = "  -1.81%  "
$ws.Range("D13").Value = "'4.357"
$ws.Range("E13").Value = "  -1.74%  "
$ws.Range("D14").Value = "1.864.95"
$ws.Range("E14").Value = "  -1.75%  "
$ws.Range("D15").Value = "'0.5400"
$ws.Range("E15").Value = "  -3.73%  "
$ws.Range("D16").Value = "0.0₅7952"
$ws.Range("D17").Value = "'64.96"
$ws.Range("E17").Value = "  -2.00%  "
$ws.Range("D18").Value = "25.992.30"
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("D19").Value = "'1.003"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D20").Value = "'4.640"
$ws.Range("E20").Value = "  -3.11%  "
$ws.Range("D21").Value = "'185.33"
$ws.Range("E21").Value = "  -1.05%  "
$ws.Range("D22").Value = "'10.03"
$ws.Range("E22").Value = "  -3.17%  "
$ws.Range("D23").Value = "'6.086"
$ws.Range("E23").Value = "  -1.34%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").Value = "'145.35"
$ws.Range("E25").Value = "  -1.83%  "
$ws.Range("D26").Value = "'7.349"
$ws.Range("E26").Value = "  -3.22%  "
$ws.Range("D27").Value = "'0.1191"
$ws.Range("E27").Value = "  -4.27%  "
$ws.Range("D28").Value = "'15.47"
$ws.Range("E28").Value = "  -3.07%  "
$ws.Range("D29").Value = "'1.369"
$ws.Range("E29").Value = "  +0.42%  "
$ws.Range("D30").Value = "'0.05974"
$ws.Range("E30").Value = "  -4.00%  "
$ws.Range("D31").Value = "'1.244"
$ws.Range("E31").Value = "  -2.70%  "
$ws.Range("D32").Value = "'3.353"
$ws.Range("E32").Value = "  -3.17%  "
$ws.Range("D33").Value = "'3.338"
$ws.Range("E33").Value = "  -2.70%  "
$ws.Range("D34").Value = "'1.609"
$ws.Range("E34").Value = "  -0.58%  "
$ws.Range("D35").Value = "'0.9697"
$ws.Range("E35").Value = "  -2.40%  "
$ws.Range("D36").Value = "'2.381"
$ws.Range("E36").Value = "  -0.90%  "
$ws.Range("E37").Value = "  +0.72%  "
$ws.Range("D38").Value = "'0.5834"
$ws.Range("E38").Value = "  -3.59%  "
$ws.Range("D39").Value = "'0.01592"
$ws.Range("E39").Value = "  -1.28%  "
$ws.Range("D40").Value = "1.047.39"
$ws.Range("E40").Value = "  -2.19%  "
$ws.Range("D41").Value = "'5.761"
$ws.Range("E41").Value = "  -6.03%  "
$ws.Range("D42").Value = "'1.003"
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("D43").Value = "'0.8406"
$ws.Range("E43").Value = "  -2.89%  "
$ws.Range("D44").Value = "'99.73"
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").Value = "1.791.93"
$ws.Range("E46").Value = "  -1.97%  "
$ws.Range("D47").Value = "'1.001"
$ws.Range("E47").Value = "  -0.24%  "
$ws.Range("D48").Value = "'54.15"
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'7.950"
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.05195"
$ws.Range("E50").Value = "  -0.95%  "
$ws.Range("D51").Value = "'0.4233"
$ws.Range("E51").Value = "  -0.50%  "
